$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet1 content updates (schema table tweaks: PK/FK columns, renamed fields,
# new notes rows at the bottom, old phone-number rows replaced).
# ---------------------------------------------------------------------------
$ws1.Range("C2").Value = "PK"
$ws1.Range("K2").Value = "Seq#"
$ws1.Range("G3").Value = "FK-->Product.ID"
$ws1.Range("L4").Value = "FK-->Product.ID"
$ws1.Range("K6").Value = "Rating (1-5)"
$ws1.Range("L6").Value = "number"
$ws1.Range("F7").Value = "User_Contact Number"
$ws1.Range("A10").Value = "Product_Contact Number"
$ws1.Range("B10").Value = "text"

# Old Landline/Mobile rows (A11:A14 + B14) are gone.
$ws1.Range("A11").ClearContents()
$ws1.Range("A12").ClearContents()
$ws1.Range("A13").ClearContents()
$ws1.Range("A14").ClearContents()
$ws1.Range("B14").ClearContents()
$ws1.Range("B15").ClearContents()

# The "Category/.../Videos&Photos" note moved from F18 up to F12.
$ws1.Range("F18").ClearContents()
$ws1.Range("F12").Value = "Category/sub-category/userid/ID/Videos&Photos"

# New freeform notes replacing the old "Votes"/"comments" rows.
$ws1.Range("A15").Value = "Once user login to website from Facebook login, Next time when he/she logs in to the facebook, we need to show the advertisement of our site with top rating products."
$ws1.Range("A16").Value = "Implement Social login plugin"
$ws1.Range("A17").Value = "use LAMP"
$ws1.Range("A18").Value = "No Payment Gateway"

# New column K best-fits to the widest header it now holds.
$ws1.Columns("K").ColumnWidth = 9.5

# ---------------------------------------------------------------------------
# Add Sheet2 (after Sheet1) with the Main_Class/Method + FrontEnd tables.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)

$ws2.Range("C4").Value = "Main_Class"
$ws2.Range("F4").Value = "FrontEnd"
$ws2.Range("C5").Value = "Method"
$ws2.Range("D5").Value = "Search_Product"
$ws2.Range("F5").Value = "Show Product based on location(maps)"
$ws2.Range("D6").Value = "Add_Product"
$ws2.Range("D7").Value = "Remove_Product"
$ws2.Range("D8").Value = "Update_Product"

$ws2.Columns("D").ColumnWidth = 14.25

# ---------------------------------------------------------------------------
# Selections matching the saved view state: Sheet1 shows L3 selected (no
# longer the active tab), Sheet2 is the active tab with F6 selected.
# ---------------------------------------------------------------------------
$ws1.Range("L3").Select()
$ws2.Range("F6").Select()
$ws2.Activate()
